# Auto-generated Excel COM-interop script applying value updates
# as described by the unified OOXML diff (Malboro_Profits workbook).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 26.25  # H6: 23.923077 -> 26.25
$ws.Cells.Item(6, 9).Value = 26.25  # I6: 23.923077 -> 26.25
$ws.Cells.Item(6, 11).Value = 78.75  # K6: 71.76923099999999 -> 78.75
$ws.Cells.Item(6, 13).Value = 33.25  # M6: 40.23076900000001 -> 33.25
$ws.Cells.Item(7, 8).Value = 2169.3333  # H7: 879.5 -> 2169.3333
$ws.Cells.Item(7, 9).Value = 8  # I7: 0 -> 8
$ws.Cells.Item(7, 10).Value = 3250  # J7: 879.5 -> 3250
$ws.Cells.Item(7, 11).Value = 8  # K7: 0 -> 8
$ws.Cells.Item(7, 12).Value = 3250  # L7: 879.5 -> 3250
$ws.Cells.Item(7, 13).Value = 104  # M7: None -> 104
$ws.Cells.Item(7, 14).Value = -3474  # N7: -1103.5 -> -3474
$ws.Cells.Item(9, 8).Value = 1362.9474  # H9: 1510.8823 -> 1362.9474
$ws.Cells.Item(9, 9).Value = 1355.5454  # I9: 1633.3334 -> 1355.5454
$ws.Cells.Item(9, 11).Value = 1355.5454  # K9: 1633.3334 -> 1355.5454
$ws.Cells.Item(9, 13).Value = -1186.5454  # M9: -1464.3334 -> -1186.5454
$ws.Cells.Item(11, 8).Value = 125.916664  # H11: 89.588234 -> 125.916664
$ws.Cells.Item(11, 9).Value = 125.916664  # I11: 89.588234 -> 125.916664
$ws.Cells.Item(11, 11).Value = 125.916664  # K11: 89.588234 -> 125.916664
$ws.Cells.Item(11, 13).Value = 14.083336  # M11: 50.411766 -> 14.083336
$ws.Cells.Item(14, 8).Value = 2169.3333  # H14: 879.5 -> 2169.3333
$ws.Cells.Item(14, 9).Value = 8  # I14: 0 -> 8
$ws.Cells.Item(14, 10).Value = 3250  # J14: 879.5 -> 3250
$ws.Cells.Item(14, 11).Value = 8  # K14: 0 -> 8
$ws.Cells.Item(14, 12).Value = 3250  # L14: 879.5 -> 3250
$ws.Cells.Item(14, 13).Value = 183  # M14: None -> 183
$ws.Cells.Item(14, 14).Value = -3632  # N14: -1261.5 -> -3632
$ws.Cells.Item(40, 8).Value = 4825.7144  # H40: 2588.3333 -> 4825.7144
$ws.Cells.Item(40, 9).Value = 4000  # I40: 3437.5 -> 4000
$ws.Cells.Item(40, 10).Value = 5156  # J40: 890 -> 5156
$ws.Cells.Item(40, 11).Value = 4000  # K40: 3437.5 -> 4000
$ws.Cells.Item(40, 12).Value = 5156  # L40: 890 -> 5156
$ws.Cells.Item(40, 13).Value = -3825  # M40: -3262.5 -> -3825
$ws.Cells.Item(40, 14).Value = -5506  # N40: -1240 -> -5506
$ws.Cells.Item(51, 8).Value = 5380  # H51: 7875 -> 5380
$ws.Cells.Item(51, 9).Value = 2500  # I51: 3600 -> 2500
$ws.Cells.Item(51, 10).Value = 7300  # J51: 9300 -> 7300
$ws.Cells.Item(51, 11).Value = 2500  # K51: 3600 -> 2500
$ws.Cells.Item(51, 12).Value = 7300  # L51: 9300 -> 7300
$ws.Cells.Item(51, 13).Value = -2016  # M51: -3116 -> -2016
$ws.Cells.Item(51, 14).Value = -8268  # N51: -10268 -> -8268
$ws.Cells.Item(62, 8).Value = 8146.9165  # H62: 8156.1665 -> 8146.9165
$ws.Cells.Item(62, 9).Value = 7386.2856  # I62: 7402.143 -> 7386.2856
$ws.Cells.Item(62, 11).Value = 7386.2856  # K62: 7402.143 -> 7386.2856
$ws.Cells.Item(62, 13).Value = -6762.2856  # M62: -6778.143 -> -6762.2856
$ws.Cells.Item(65, 8).Value = 8146.9165  # H65: 8156.1665 -> 8146.9165
$ws.Cells.Item(65, 9).Value = 7386.2856  # I65: 7402.143 -> 7386.2856
$ws.Cells.Item(65, 11).Value = 36931.428  # K65: 37010.715 -> 36931.428
$ws.Cells.Item(65, 13).Value = -33811.428  # M65: -33890.715 -> -33811.428
$ws.Cells.Item(76, 8).Value = 7749  # H76: 7999 -> 7749
$ws.Cells.Item(76, 9).Value = 7749  # I76: 7999 -> 7749
$ws.Cells.Item(76, 11).Value = 7749  # K76: 7999 -> 7749
$ws.Cells.Item(76, 13).Value = -7434  # M76: -7684 -> -7434
$ws.Cells.Item(79, 8).Value = 7749  # H79: 7999 -> 7749
$ws.Cells.Item(79, 9).Value = 7749  # I79: 7999 -> 7749
$ws.Cells.Item(79, 11).Value = 7749  # K79: 7999 -> 7749
$ws.Cells.Item(79, 13).Value = -6657  # M79: -6907 -> -6657
$ws.Cells.Item(97, 8).Value = 757.2857  # H97: 902.4 -> 757.2857
$ws.Cells.Item(97, 10).Value = 833.5  # J97: 1053 -> 833.5
$ws.Cells.Item(97, 12).Value = 2500.5  # L97: 3159 -> 2500.5
$ws.Cells.Item(97, 14).Value = -3492.5  # N97: -4151 -> -3492.5
$ws.Cells.Item(98, 9).Value = 11376.25  # I98: 9501 -> 11376.25
$ws.Cells.Item(98, 10).Value = 5000  # J98: 6000 -> 5000
$ws.Cells.Item(98, 11).Value = 11376.25  # K98: 9501 -> 11376.25
$ws.Cells.Item(98, 12).Value = 5000  # L98: 6000 -> 5000
$ws.Cells.Item(98, 13).Value = -9878.25  # M98: -8003 -> -9878.25
$ws.Cells.Item(98, 14).Value = -7996  # N98: -8996 -> -7996
$ws.Cells.Item(122, 9).Value = 11376.25  # I122: 9501 -> 11376.25
$ws.Cells.Item(122, 10).Value = 5000  # J122: 6000 -> 5000
$ws.Cells.Item(122, 11).Value = 34128.75  # K122: 28503 -> 34128.75
$ws.Cells.Item(122, 12).Value = 15000  # L122: 18000 -> 15000
$ws.Cells.Item(122, 13).Value = -31678.75  # M122: -26053 -> -31678.75
$ws.Cells.Item(122, 14).Value = -19900  # N122: -22900 -> -19900
$ws.Cells.Item(129, 8).Value = 1628.5  # H129: 765.6667 -> 1628.5
$ws.Cells.Item(129, 10).Value = 4217  # J129: 0 -> 4217
$ws.Cells.Item(129, 12).Value = 12651  # L129: 0 -> 12651
$ws.Cells.Item(129, 14).Value = -22651  # N129: None -> -22651
$ws.Cells.Item(131, 8).Value = 5535  # H131: 7052.5 -> 5535
$ws.Cells.Item(131, 9).Value = 2500  # I131: 0 -> 2500
$ws.Cells.Item(131, 11).Value = 7500  # K131: 0 -> 7500
$ws.Cells.Item(131, 13).Value = -2460  # M131: None -> -2460
$ws.Cells.Item(132, 8).Value = 10755.758  # H132: 11004.375 -> 10755.758
$ws.Cells.Item(132, 9).Value = 9361.959999999999  # I132: 9369.959999999999 -> 9361.959999999999
$ws.Cells.Item(132, 10).Value = 15111.375  # J132: 16841.572 -> 15111.375
$ws.Cells.Item(132, 11).Value = 28085.88  # K132: 28109.88 -> 28085.88
$ws.Cells.Item(132, 12).Value = 45334.125  # L132: 50524.716 -> 45334.125
$ws.Cells.Item(132, 13).Value = -25555.88  # M132: -25579.88 -> -25555.88
$ws.Cells.Item(132, 14).Value = -50394.125  # N132: -55584.716 -> -50394.125
$ws.Cells.Item(135, 8).Value = 1593.9565  # H135: 1379.1482 -> 1593.9565
$ws.Cells.Item(135, 9).Value = 1667.8572  # I135: 1473.7916 -> 1667.8572
$ws.Cells.Item(135, 10).Value = 818  # J135: 622 -> 818
$ws.Cells.Item(135, 11).Value = 15010.7148  # K135: 13264.1244 -> 15010.7148
$ws.Cells.Item(135, 12).Value = 7362  # L135: 5598 -> 7362
$ws.Cells.Item(135, 13).Value = -12475.7148  # M135: -10729.1244 -> -12475.7148
$ws.Cells.Item(135, 14).Value = -12432  # N135: -10668 -> -12432
$ws.Cells.Item(141, 8).Value = 2710.125  # H141: 2240.1428 -> 2710.125
$ws.Cells.Item(141, 10).Value = 6299.6665  # J141: 6449.5 -> 6299.6665
$ws.Cells.Item(141, 12).Value = 18898.9995  # L141: 19348.5 -> 18898.9995
$ws.Cells.Item(141, 14).Value = -29258.9995  # N141: -29708.5 -> -29258.9995
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1041.4054  # H2: 1053.1351 -> 1041.4054
$ws.Cells.Item(2, 9).Value = 774.76  # I2: 792.12 -> 774.76
$ws.Cells.Item(2, 11).Value = 774.76  # K2: 792.12 -> 774.76
$ws.Cells.Item(2, 13).Value = -661.76  # M2: -679.12 -> -661.76
$ws.Cells.Item(8, 8).Value = 2060.8333  # H8: 3666.6667 -> 2060.8333
$ws.Cells.Item(8, 9).Value = 1749.75  # I8: 3000 -> 1749.75
$ws.Cells.Item(8, 10).Value = 2683  # J8: 5000 -> 2683
$ws.Cells.Item(8, 11).Value = 1749.75  # K8: 3000 -> 1749.75
$ws.Cells.Item(8, 12).Value = 2683  # L8: 5000 -> 2683
$ws.Cells.Item(8, 13).Value = -1605.75  # M8: -2856 -> -1605.75
$ws.Cells.Item(8, 14).Value = -2971  # N8: -5288 -> -2971
$ws.Cells.Item(16, 8).Value = 732.6  # H16: 645.375 -> 732.6
$ws.Cells.Item(16, 9).Value = 732.6  # I16: 645.375 -> 732.6
$ws.Cells.Item(16, 11).Value = 732.6  # K16: 645.375 -> 732.6
$ws.Cells.Item(16, 13).Value = -445.6  # M16: -358.375 -> -445.6
$ws.Cells.Item(25, 8).Value = 5750  # H25: 5085.7144 -> 5750
$ws.Cells.Item(25, 9).Value = 3046.4  # I25: 2966.4 -> 3046.4
$ws.Cells.Item(25, 10).Value = 10256  # J25: 10384 -> 10256
$ws.Cells.Item(25, 11).Value = 3046.4  # K25: 2966.4 -> 3046.4
$ws.Cells.Item(25, 12).Value = 10256  # L25: 10384 -> 10256
$ws.Cells.Item(25, 13).Value = -2644.4  # M25: -2564.4 -> -2644.4
$ws.Cells.Item(25, 14).Value = -11060  # N25: -11188 -> -11060
$ws.Cells.Item(32, 8).Value = 10959.73  # H32: 11372.2 -> 10959.73
$ws.Cells.Item(32, 9).Value = 8407.137000000001  # I32: 8776.619000000001 -> 8407.137000000001
$ws.Cells.Item(32, 11).Value = 8407.137000000001  # K32: 8776.619000000001 -> 8407.137000000001
$ws.Cells.Item(32, 13).Value = -8120.137000000001  # M32: -8489.619000000001 -> -8120.137000000001
$ws.Cells.Item(34, 8).Value = 5000  # H34: 40000 -> 5000
$ws.Cells.Item(34, 10).Value = 5000  # J34: 40000 -> 5000
$ws.Cells.Item(34, 12).Value = 5000  # L34: 40000 -> 5000
$ws.Cells.Item(34, 14).Value = -5542  # N34: -40542 -> -5542
$ws.Cells.Item(45, 8).Value = 2045.55  # H45: 2053.0952 -> 2045.55
$ws.Cells.Item(45, 9).Value = 2102.7856  # I45: 2110.5 -> 2102.7856
$ws.Cells.Item(45, 10).Value = 1912  # J45: 1938.2858 -> 1912
$ws.Cells.Item(45, 11).Value = 2102.7856  # K45: 2110.5 -> 2102.7856
$ws.Cells.Item(45, 12).Value = 1912  # L45: 1938.2858 -> 1912
$ws.Cells.Item(45, 13).Value = -1725.7856  # M45: -1733.5 -> -1725.7856
$ws.Cells.Item(45, 14).Value = -2666  # N45: -2692.2858 -> -2666
$ws.Cells.Item(61, 8).Value = 14019.939  # H61: 14783.29 -> 14019.939
$ws.Cells.Item(61, 9).Value = 3953.4443  # I61: 4653.3335 -> 3953.4443
$ws.Cells.Item(61, 10).Value = 17794.875  # J61: 17214.48 -> 17794.875
$ws.Cells.Item(61, 11).Value = 3953.4443  # K61: 4653.3335 -> 3953.4443
$ws.Cells.Item(61, 12).Value = 17794.875  # L61: 17214.48 -> 17794.875
$ws.Cells.Item(61, 13).Value = -3741.4443  # M61: -4441.3335 -> -3741.4443
$ws.Cells.Item(61, 14).Value = -18218.875  # N61: -17638.48 -> -18218.875
$ws.Cells.Item(74, 8).Value = 17317.346  # H74: 17982.76 -> 17317.346
$ws.Cells.Item(74, 9).Value = 1487.125  # I74: 1602.1428 -> 1487.125
$ws.Cells.Item(74, 11).Value = 1487.125  # K74: 1602.1428 -> 1487.125
$ws.Cells.Item(74, 13).Value = -613.125  # M74: -728.1428000000001 -> -613.125
$ws.Cells.Item(77, 8).Value = 17317.346  # H77: 17982.76 -> 17317.346
$ws.Cells.Item(77, 9).Value = 1487.125  # I77: 1602.1428 -> 1487.125
$ws.Cells.Item(77, 11).Value = 7435.625  # K77: 8010.714 -> 7435.625
$ws.Cells.Item(77, 13).Value = -3067.625  # M77: -3642.714 -> -3067.625
$ws.Cells.Item(102, 8).Value = 15300.235  # H102: 20539.916 -> 15300.235
$ws.Cells.Item(102, 9).Value = 3469.6155  # I102: 3935 -> 3469.6155
$ws.Cells.Item(102, 11).Value = 3469.6155  # K102: 3935 -> 3469.6155
$ws.Cells.Item(102, 13).Value = -1847.6155  # M102: -2313 -> -1847.6155
$ws.Cells.Item(116, 8).Value = 1041.4054  # H116: 1053.1351 -> 1041.4054
$ws.Cells.Item(116, 9).Value = 774.76  # I116: 792.12 -> 774.76
$ws.Cells.Item(116, 11).Value = 774.76  # K116: 792.12 -> 774.76
$ws.Cells.Item(116, 13).Value = 1519.24  # M116: 1501.88 -> 1519.24
$ws.Cells.Item(122, 8).Value = 3615.4  # H122: 4522 -> 3615.4
$ws.Cells.Item(122, 9).Value = 2021  # I122: 2241.375 -> 2021
$ws.Cells.Item(122, 10).Value = 8000  # J122: 7128.4287 -> 8000
$ws.Cells.Item(122, 11).Value = 6063  # K122: 6724.125 -> 6063
$ws.Cells.Item(122, 12).Value = 24000  # L122: 21385.2861 -> 24000
$ws.Cells.Item(122, 13).Value = -3613  # M122: -4274.125 -> -3613
$ws.Cells.Item(122, 14).Value = -28900  # N122: -26285.2861 -> -28900
$ws.Cells.Item(126, 8).Value = 8311.875  # H126: 7713.5713 -> 8311.875
$ws.Cells.Item(126, 9).Value = 8311.875  # I126: 7713.5713 -> 8311.875
$ws.Cells.Item(126, 11).Value = 24935.625  # K126: 23140.7139 -> 24935.625
$ws.Cells.Item(126, 13).Value = -22465.625  # M126: -20670.7139 -> -22465.625
$ws.Cells.Item(130, 8).Value = 54677.125  # H130: 52600.777 -> 54677.125
$ws.Cells.Item(130, 10).Value = 54677.125  # J130: 52600.777 -> 54677.125
$ws.Cells.Item(130, 12).Value = 54677.125  # L130: 52600.777 -> 54677.125
$ws.Cells.Item(130, 14).Value = -64717.125  # N130: -62640.777 -> -64717.125
$ws.Cells.Item(132, 8).Value = 2336336.5  # H132: 2575861 -> 2336336.5
$ws.Cells.Item(132, 9).Value = 3657.5  # I132: 4105.0835 -> 3657.5
$ws.Cells.Item(132, 11).Value = 10972.5  # K132: 12315.2505 -> 10972.5
$ws.Cells.Item(132, 13).Value = -8442.5  # M132: -9785.250499999998 -> -8442.5
$ws.Cells.Item(136, 8).Value = 14019.939  # H136: 14783.29 -> 14019.939
$ws.Cells.Item(136, 9).Value = 3953.4443  # I136: 4653.3335 -> 3953.4443
$ws.Cells.Item(136, 10).Value = 17794.875  # J136: 17214.48 -> 17794.875
$ws.Cells.Item(136, 11).Value = 11860.3329  # K136: 13960.0005 -> 11860.3329
$ws.Cells.Item(136, 12).Value = 53384.625  # L136: 51643.44 -> 53384.625
$ws.Cells.Item(136, 13).Value = -9310.332900000001  # M136: -11410.0005 -> -9310.332900000001
$ws.Cells.Item(136, 14).Value = -58484.625  # N136: -56743.44 -> -58484.625
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1041.4054  # H3: 1053.1351 -> 1041.4054
$ws.Cells.Item(3, 9).Value = 774.76  # I3: 792.12 -> 774.76
$ws.Cells.Item(3, 11).Value = 774.76  # K3: 792.12 -> 774.76
$ws.Cells.Item(3, 13).Value = -660.76  # M3: -678.12 -> -660.76
$ws.Cells.Item(20, 8).Value = 18746.979  # H20: 19130.467 -> 18746.979
$ws.Cells.Item(20, 9).Value = 6191  # I20: 6174.478 -> 6191
$ws.Cells.Item(20, 10).Value = 31302.957  # J20: 32675.363 -> 31302.957
$ws.Cells.Item(20, 11).Value = 6191  # K20: 6174.478 -> 6191
$ws.Cells.Item(20, 12).Value = 31302.957  # L20: 32675.363 -> 31302.957
$ws.Cells.Item(20, 13).Value = -5944  # M20: -5927.478 -> -5944
$ws.Cells.Item(20, 14).Value = -31796.957  # N20: -33169.363 -> -31796.957
$ws.Cells.Item(31, 8).Value = 6000  # H31: 0 -> 6000
$ws.Cells.Item(31, 10).Value = 6000  # J31: 0 -> 6000
$ws.Cells.Item(31, 12).Value = 6000  # L31: 0 -> 6000
$ws.Cells.Item(31, 14).Value = -6504  # N31: None -> -6504
$ws.Cells.Item(128, 8).Value = 7626.273  # H128: 7680.8184 -> 7626.273
$ws.Cells.Item(128, 9).Value = 7626.273  # I128: 7680.8184 -> 7626.273
$ws.Cells.Item(128, 11).Value = 22878.819  # K128: 23042.4552 -> 22878.819
$ws.Cells.Item(128, 13).Value = -20388.819  # M128: -20552.4552 -> -20388.819
$ws.Cells.Item(132, 8).Value = 78915.57000000001  # H132: 78801.25 -> 78915.57000000001
$ws.Cells.Item(132, 10).Value = 78915.57000000001  # J132: 78801.25 -> 78915.57000000001
$ws.Cells.Item(132, 12).Value = 78915.57000000001  # L132: 78801.25 -> 78915.57000000001
$ws.Cells.Item(132, 14).Value = -89035.57000000001  # N132: -88921.25 -> -89035.57000000001
$ws.Cells.Item(134, 8).Value = 10495  # H134: 11196.64 -> 10495
$ws.Cells.Item(134, 9).Value = 2239.9443  # I134: 2268.8235 -> 2239.9443
$ws.Cells.Item(134, 10).Value = 27005.111  # J134: 30168.25 -> 27005.111
$ws.Cells.Item(134, 11).Value = 6719.8329  # K134: 6806.470499999999 -> 6719.8329
$ws.Cells.Item(134, 12).Value = 81015.333  # L134: 90504.75 -> 81015.333
$ws.Cells.Item(134, 13).Value = -4184.8329  # M134: -4271.470499999999 -> -4184.8329
$ws.Cells.Item(134, 14).Value = -86085.333  # N134: -95574.75 -> -86085.333
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 167.10527  # H7: 148.22728 -> 167.10527
$ws.Cells.Item(7, 9).Value = 94.30768999999999  # I7: 85.066666 -> 94.30768999999999
$ws.Cells.Item(7, 10).Value = 324.83334  # J7: 283.57144 -> 324.83334
$ws.Cells.Item(7, 11).Value = 94.30768999999999  # K7: 85.066666 -> 94.30768999999999
$ws.Cells.Item(7, 12).Value = 324.83334  # L7: 283.57144 -> 324.83334
$ws.Cells.Item(7, 13).Value = 18.69231000000001  # M7: 27.933334 -> 18.69231000000001
$ws.Cells.Item(7, 14).Value = -550.83334  # N7: -509.57144 -> -550.83334
$ws.Cells.Item(16, 8).Value = 1507.25  # H16: 1534.7333 -> 1507.25
$ws.Cells.Item(16, 10).Value = 1840.8334  # J16: 1990 -> 1840.8334
$ws.Cells.Item(16, 12).Value = 1840.8334  # L16: 1990 -> 1840.8334
$ws.Cells.Item(16, 14).Value = -2414.8334  # N16: -2564 -> -2414.8334
$ws.Cells.Item(22, 8).Value = 1556.3334  # H22: 1396.8889 -> 1556.3334
$ws.Cells.Item(22, 9).Value = 1049  # I22: 827.8570999999999 -> 1049
$ws.Cells.Item(22, 10).Value = 1740.8182  # J22: 1759 -> 1740.8182
$ws.Cells.Item(22, 11).Value = 1049  # K22: 827.8570999999999 -> 1049
$ws.Cells.Item(22, 12).Value = 1740.8182  # L22: 1759 -> 1740.8182
$ws.Cells.Item(22, 13).Value = -699  # M22: -477.8570999999999 -> -699
$ws.Cells.Item(22, 14).Value = -2440.8182  # N22: -2459 -> -2440.8182
$ws.Cells.Item(45, 8).Value = 1000  # H45: 0 -> 1000
$ws.Cells.Item(45, 9).Value = 1000  # I45: 0 -> 1000
$ws.Cells.Item(45, 11).Value = 1000  # K45: 0 -> 1000
$ws.Cells.Item(45, 13).Value = -407  # M45: None -> -407
$ws.Cells.Item(58, 8).Value = 27538.375  # H58: 29011.732 -> 27538.375
$ws.Cells.Item(58, 9).Value = 15037.5  # I58: 19834.25 -> 15037.5
$ws.Cells.Item(58, 10).Value = 35038.9  # J58: 32349 -> 35038.9
$ws.Cells.Item(58, 11).Value = 15037.5  # K58: 19834.25 -> 15037.5
$ws.Cells.Item(58, 12).Value = 35038.9  # L58: 32349 -> 35038.9
$ws.Cells.Item(58, 13).Value = -14834.5  # M58: -19631.25 -> -14834.5
$ws.Cells.Item(58, 14).Value = -35444.9  # N58: -32755 -> -35444.9
$ws.Cells.Item(70, 8).Value = 17000  # H70: 16666.666 -> 17000
$ws.Cells.Item(70, 10).Value = 17000  # J70: 16666.666 -> 17000
$ws.Cells.Item(70, 12).Value = 17000  # L70: 16666.666 -> 17000
$ws.Cells.Item(70, 14).Value = -17630  # N70: -17296.666 -> -17630
$ws.Cells.Item(73, 8).Value = 17000  # H73: 16666.666 -> 17000
$ws.Cells.Item(73, 10).Value = 17000  # J73: 16666.666 -> 17000
$ws.Cells.Item(73, 12).Value = 17000  # L73: 16666.666 -> 17000
$ws.Cells.Item(73, 14).Value = -19184  # N73: -18850.666 -> -19184
$ws.Cells.Item(113, 8).Value = 1507.25  # H113: 1534.7333 -> 1507.25
$ws.Cells.Item(113, 10).Value = 1840.8334  # J113: 1990 -> 1840.8334
$ws.Cells.Item(113, 12).Value = 1840.8334  # L113: 1990 -> 1840.8334
$ws.Cells.Item(113, 14).Value = -6180.8334  # N113: -6330 -> -6180.8334
$ws.Cells.Item(127, 8).Value = 20000  # H127: 0 -> 20000
$ws.Cells.Item(127, 10).Value = 20000  # J127: 0 -> 20000
$ws.Cells.Item(127, 12).Value = 20000  # L127: 0 -> 20000
$ws.Cells.Item(127, 14).Value = -29920  # N127: None -> -29920
$ws.Cells.Item(136, 8).Value = 27538.375  # H136: 29011.732 -> 27538.375
$ws.Cells.Item(136, 9).Value = 15037.5  # I136: 19834.25 -> 15037.5
$ws.Cells.Item(136, 10).Value = 35038.9  # J136: 32349 -> 35038.9
$ws.Cells.Item(136, 11).Value = 45112.5  # K136: 59502.75 -> 45112.5
$ws.Cells.Item(136, 12).Value = 105116.7  # L136: 97047 -> 105116.7
$ws.Cells.Item(136, 13).Value = -42562.5  # M136: -56952.75 -> -42562.5
$ws.Cells.Item(136, 14).Value = -110216.7  # N136: -102147 -> -110216.7
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 276.6316  # H2: 265.35 -> 276.6316
$ws.Cells.Item(2, 9).Value = 300.7647  # I2: 286.8889 -> 300.7647
$ws.Cells.Item(2, 11).Value = 1804.5882  # K2: 1721.3334 -> 1804.5882
$ws.Cells.Item(2, 13).Value = -1691.5882  # M2: -1608.3334 -> -1691.5882
$ws.Cells.Item(11, 8).Value = 683.0454999999999  # H11: 651.56525 -> 683.0454999999999
$ws.Cells.Item(11, 9).Value = 687.82355  # I11: 639.5789 -> 687.82355
$ws.Cells.Item(11, 10).Value = 666.8  # J11: 708.5 -> 666.8
$ws.Cells.Item(11, 11).Value = 2063.47065  # K11: 1918.7367 -> 2063.47065
$ws.Cells.Item(11, 12).Value = 2000.4  # L11: 2125.5 -> 2000.4
$ws.Cells.Item(11, 13).Value = -1923.47065  # M11: -1778.7367 -> -1923.47065
$ws.Cells.Item(11, 14).Value = -2280.4  # N11: -2405.5 -> -2280.4
$ws.Cells.Item(26, 8).Value = 2061  # H26: 2485.4 -> 2061
$ws.Cells.Item(26, 10).Value = 849.75  # J26: 699.5 -> 849.75
$ws.Cells.Item(26, 12).Value = 2549.25  # L26: 2098.5 -> 2549.25
$ws.Cells.Item(26, 14).Value = -3125.25  # N26: -2674.5 -> -3125.25
$ws.Cells.Item(34, 8).Value = 1054.25  # H34: 1260.2 -> 1054.25
$ws.Cells.Item(34, 9).Value = 1269.5555  # I34: 1425.25 -> 1269.5555
$ws.Cells.Item(34, 10).Value = 408.33334  # J34: 600 -> 408.33334
$ws.Cells.Item(34, 11).Value = 3808.6665  # K34: 4275.75 -> 3808.6665
$ws.Cells.Item(34, 12).Value = 1225.00002  # L34: 1800 -> 1225.00002
$ws.Cells.Item(34, 13).Value = -3724.6665  # M34: -4191.75 -> -3724.6665
$ws.Cells.Item(34, 14).Value = -1393.00002  # N34: -1968 -> -1393.00002
$ws.Cells.Item(64, 8).Value = 3999.75  # H64: 2000 -> 3999.75
$ws.Cells.Item(64, 9).Value = 1999.6666  # I64: 2000 -> 1999.6666
$ws.Cells.Item(64, 10).Value = 10000  # J64: 0 -> 10000
$ws.Cells.Item(64, 11).Value = 5998.9998  # K64: 6000 -> 5998.9998
$ws.Cells.Item(64, 12).Value = 30000  # L64: 0 -> 30000
$ws.Cells.Item(64, 13).Value = -5728.9998  # M64: -5730 -> -5728.9998
$ws.Cells.Item(64, 14).Value = -30540  # N64: None -> -30540
$ws.Cells.Item(67, 8).Value = 3999.75  # H67: 2000 -> 3999.75
$ws.Cells.Item(67, 9).Value = 1999.6666  # I67: 2000 -> 1999.6666
$ws.Cells.Item(67, 10).Value = 10000  # J67: 0 -> 10000
$ws.Cells.Item(67, 11).Value = 5998.9998  # K67: 6000 -> 5998.9998
$ws.Cells.Item(67, 12).Value = 30000  # L67: 0 -> 30000
$ws.Cells.Item(67, 13).Value = -5062.9998  # M67: -5064 -> -5062.9998
$ws.Cells.Item(67, 14).Value = -31872  # N67: None -> -31872
$ws.Cells.Item(68, 8).Value = 1687  # H68: 1687.125 -> 1687
$ws.Cells.Item(68, 10).Value = 1856.5714  # J68: 1856.7142 -> 1856.5714
$ws.Cells.Item(68, 12).Value = 5569.7142  # L68: 5570.142599999999 -> 5569.7142
$ws.Cells.Item(68, 14).Value = -7191.7142  # N68: -7192.142599999999 -> -7191.7142
$ws.Cells.Item(69, 8).Value = 12666.333  # H69: 9000 -> 12666.333
$ws.Cells.Item(69, 9).Value = 8999  # I69: 0 -> 8999
$ws.Cells.Item(69, 10).Value = 14500  # J69: 9000 -> 14500
$ws.Cells.Item(69, 11).Value = 26997  # K69: 0 -> 26997
$ws.Cells.Item(69, 12).Value = 43500  # L69: 27000 -> 43500
$ws.Cells.Item(69, 13).Value = -26186  # M69: None -> -26186
$ws.Cells.Item(69, 14).Value = -45122  # N69: -28622 -> -45122
$ws.Cells.Item(70, 8).Value = 10999.5  # H70: 13499.667 -> 10999.5
$ws.Cells.Item(70, 9).Value = 1000  # I70: 0 -> 1000
$ws.Cells.Item(70, 10).Value = 14332.667  # J70: 13499.667 -> 14332.667
$ws.Cells.Item(70, 11).Value = 3000  # K70: 0 -> 3000
$ws.Cells.Item(70, 12).Value = 42998.001  # L70: 40499.001 -> 42998.001
$ws.Cells.Item(70, 13).Value = -2685  # M70: None -> -2685
$ws.Cells.Item(70, 14).Value = -43628.001  # N70: -41129.001 -> -43628.001
$ws.Cells.Item(71, 8).Value = 1687  # H71: 1687.125 -> 1687
$ws.Cells.Item(71, 10).Value = 1856.5714  # J71: 1856.7142 -> 1856.5714
$ws.Cells.Item(71, 12).Value = 16709.1426  # L71: 16710.4278 -> 16709.1426
$ws.Cells.Item(71, 14).Value = -24821.1426  # N71: -24822.4278 -> -24821.1426
$ws.Cells.Item(72, 8).Value = 12666.333  # H72: 9000 -> 12666.333
$ws.Cells.Item(72, 9).Value = 8999  # I72: 0 -> 8999
$ws.Cells.Item(72, 10).Value = 14500  # J72: 9000 -> 14500
$ws.Cells.Item(72, 11).Value = 80991  # K72: 0 -> 80991
$ws.Cells.Item(72, 12).Value = 130500  # L72: 81000 -> 130500
$ws.Cells.Item(72, 13).Value = -76935  # M72: None -> -76935
$ws.Cells.Item(72, 14).Value = -138612  # N72: -89112 -> -138612
$ws.Cells.Item(73, 8).Value = 10999.5  # H73: 13499.667 -> 10999.5
$ws.Cells.Item(73, 9).Value = 1000  # I73: 0 -> 1000
$ws.Cells.Item(73, 10).Value = 14332.667  # J73: 13499.667 -> 14332.667
$ws.Cells.Item(73, 11).Value = 3000  # K73: 0 -> 3000
$ws.Cells.Item(73, 12).Value = 42998.001  # L73: 40499.001 -> 42998.001
$ws.Cells.Item(73, 13).Value = -1908  # M73: None -> -1908
$ws.Cells.Item(73, 14).Value = -45182.001  # N73: -42683.001 -> -45182.001
$ws.Cells.Item(74, 8).Value = 17500  # H74: 21250 -> 17500
$ws.Cells.Item(74, 10).Value = 17500  # J74: 21250 -> 17500
$ws.Cells.Item(74, 12).Value = 52500  # L74: 63750 -> 52500
$ws.Cells.Item(74, 14).Value = -54622  # N74: -65872 -> -54622
$ws.Cells.Item(76, 8).Value = 12271  # H76: 5438 -> 12271
$ws.Cells.Item(76, 9).Value = 6813  # I76: 5438 -> 6813
$ws.Cells.Item(76, 10).Value = 15000  # J76: 0 -> 15000
$ws.Cells.Item(76, 11).Value = 20439  # K76: 16314 -> 20439
$ws.Cells.Item(76, 12).Value = 45000  # L76: 0 -> 45000
$ws.Cells.Item(76, 13).Value = -20056  # M76: -15931 -> -20056
$ws.Cells.Item(76, 14).Value = -45766  # N76: None -> -45766
$ws.Cells.Item(77, 8).Value = 17500  # H77: 21250 -> 17500
$ws.Cells.Item(77, 10).Value = 17500  # J77: 21250 -> 17500
$ws.Cells.Item(77, 12).Value = 157500  # L77: 191250 -> 157500
$ws.Cells.Item(77, 14).Value = -168108  # N77: -201858 -> -168108
$ws.Cells.Item(79, 8).Value = 12271  # H79: 5438 -> 12271
$ws.Cells.Item(79, 9).Value = 6813  # I79: 5438 -> 6813
$ws.Cells.Item(79, 10).Value = 15000  # J79: 0 -> 15000
$ws.Cells.Item(79, 11).Value = 20439  # K79: 16314 -> 20439
$ws.Cells.Item(79, 12).Value = 45000  # L79: 0 -> 45000
$ws.Cells.Item(79, 13).Value = -19113  # M79: -14988 -> -19113
$ws.Cells.Item(79, 14).Value = -47652  # N79: None -> -47652
$ws.Cells.Item(80, 8).Value = 11844.909  # H80: 12532.833 -> 11844.909
$ws.Cells.Item(80, 9).Value = 8627.714  # I80: 8642 -> 8627.714
$ws.Cells.Item(80, 10).Value = 17475  # J80: 17980 -> 17475
$ws.Cells.Item(80, 11).Value = 25883.142  # K80: 25926 -> 25883.142
$ws.Cells.Item(80, 12).Value = 52425  # L80: 53940 -> 52425
$ws.Cells.Item(80, 13).Value = -24947.142  # M80: -24990 -> -24947.142
$ws.Cells.Item(80, 14).Value = -54297  # N80: -55812 -> -54297
$ws.Cells.Item(81, 8).Value = 16670333  # H81: 50000000 -> 16670333
$ws.Cells.Item(81, 9).Value = 1000  # I81: 0 -> 1000
$ws.Cells.Item(81, 10).Value = 25005000  # J81: 50000000 -> 25005000
$ws.Cells.Item(81, 11).Value = 3000  # K81: 0 -> 3000
$ws.Cells.Item(81, 12).Value = 75015000  # L81: 150000000 -> 75015000
$ws.Cells.Item(81, 13).Value = -1877  # M81: None -> -1877
$ws.Cells.Item(81, 14).Value = -75017246  # N81: -150002246 -> -75017246
$ws.Cells.Item(82, 8).Value = 15000  # H82: 0 -> 15000
$ws.Cells.Item(82, 10).Value = 15000  # J82: 0 -> 15000
$ws.Cells.Item(82, 12).Value = 45000  # L82: 0 -> 45000
$ws.Cells.Item(82, 14).Value = -45812  # N82: None -> -45812
$ws.Cells.Item(83, 8).Value = 11844.909  # H83: 12532.833 -> 11844.909
$ws.Cells.Item(83, 9).Value = 8627.714  # I83: 8642 -> 8627.714
$ws.Cells.Item(83, 10).Value = 17475  # J83: 17980 -> 17475
$ws.Cells.Item(83, 11).Value = 77649.42600000001  # K83: 77778 -> 77649.42600000001
$ws.Cells.Item(83, 12).Value = 157275  # L83: 161820 -> 157275
$ws.Cells.Item(83, 13).Value = -72969.42600000001  # M83: -73098 -> -72969.42600000001
$ws.Cells.Item(83, 14).Value = -166635  # N83: -171180 -> -166635
$ws.Cells.Item(84, 8).Value = 16670333  # H84: 50000000 -> 16670333
$ws.Cells.Item(84, 9).Value = 1000  # I84: 0 -> 1000
$ws.Cells.Item(84, 10).Value = 25005000  # J84: 50000000 -> 25005000
$ws.Cells.Item(84, 11).Value = 9000  # K84: 0 -> 9000
$ws.Cells.Item(84, 12).Value = 225045000  # L84: 450000000 -> 225045000
$ws.Cells.Item(84, 13).Value = -3384  # M84: None -> -3384
$ws.Cells.Item(84, 14).Value = -225056232  # N84: -450011232 -> -225056232
$ws.Cells.Item(85, 8).Value = 15000  # H85: 0 -> 15000
$ws.Cells.Item(85, 10).Value = 15000  # J85: 0 -> 15000
$ws.Cells.Item(85, 12).Value = 45000  # L85: 0 -> 45000
$ws.Cells.Item(85, 14).Value = -47808  # N85: None -> -47808
$ws.Cells.Item(86, 8).Value = 742.38464  # H86: 740.08 -> 742.38464
$ws.Cells.Item(86, 9).Value = 757.8421  # I86: 729.95 -> 757.8421
$ws.Cells.Item(86, 10).Value = 700.4286  # J86: 780.6 -> 700.4286
$ws.Cells.Item(86, 11).Value = 2273.5263  # K86: 2189.85 -> 2273.5263
$ws.Cells.Item(86, 12).Value = 2101.2858  # L86: 2341.8 -> 2101.2858
$ws.Cells.Item(86, 13).Value = -1087.5263  # M86: -1003.85 -> -1087.5263
$ws.Cells.Item(86, 14).Value = -4473.2858  # N86: -4713.8 -> -4473.2858
$ws.Cells.Item(87, 8).Value = 14400  # H87: 10250 -> 14400
$ws.Cells.Item(87, 9).Value = 0  # I87: 4000 -> 0
$ws.Cells.Item(87, 10).Value = 14400  # J87: 16500 -> 14400
$ws.Cells.Item(87, 11).Value = 0  # K87: 12000 -> 0
$ws.Cells.Item(87, 12).Value = 43200  # L87: 49500 -> 43200
$ws.Cells.Item(87, 13).Value = $null  # M87: was -10752, now removed
$ws.Cells.Item(87, 14).Value = -45696  # N87: -51996 -> -45696
$ws.Cells.Item(88, 8).Value = 7000  # H88: 9700 -> 7000
$ws.Cells.Item(88, 10).Value = 7000  # J88: 9700 -> 7000
$ws.Cells.Item(88, 12).Value = 21000  # L88: 29100 -> 21000
$ws.Cells.Item(88, 14).Value = -21856  # N88: -29956 -> -21856
$ws.Cells.Item(89, 8).Value = 742.38464  # H89: 740.08 -> 742.38464
$ws.Cells.Item(89, 9).Value = 757.8421  # I89: 729.95 -> 757.8421
$ws.Cells.Item(89, 10).Value = 700.4286  # J89: 780.6 -> 700.4286
$ws.Cells.Item(89, 11).Value = 6820.5789  # K89: 6569.55 -> 6820.5789
$ws.Cells.Item(89, 12).Value = 6303.8574  # L89: 7025.400000000001 -> 6303.8574
$ws.Cells.Item(89, 13).Value = -892.5788999999995  # M89: -641.5500000000002 -> -892.5788999999995
$ws.Cells.Item(89, 14).Value = -18159.8574  # N89: -18881.4 -> -18159.8574
$ws.Cells.Item(90, 8).Value = 14400  # H90: 10250 -> 14400
$ws.Cells.Item(90, 9).Value = 0  # I90: 4000 -> 0
$ws.Cells.Item(90, 10).Value = 14400  # J90: 16500 -> 14400
$ws.Cells.Item(90, 11).Value = 0  # K90: 36000 -> 0
$ws.Cells.Item(90, 12).Value = 129600  # L90: 148500 -> 129600
$ws.Cells.Item(90, 13).Value = $null  # M90: was -29760, now removed
$ws.Cells.Item(90, 14).Value = -142080  # N90: -160980 -> -142080
$ws.Cells.Item(91, 8).Value = 7000  # H91: 9700 -> 7000
$ws.Cells.Item(91, 10).Value = 7000  # J91: 9700 -> 7000
$ws.Cells.Item(91, 12).Value = 21000  # L91: 29100 -> 21000
$ws.Cells.Item(91, 14).Value = -23964  # N91: -32064 -> -23964
$ws.Cells.Item(131, 8).Value = 1495.1327  # H131: 1492.65 -> 1495.1327
$ws.Cells.Item(131, 9).Value = 1350  # I131: 1333 -> 1350
$ws.Cells.Item(131, 10).Value = 1498.1562  # J131: 1497.5876 -> 1498.1562
$ws.Cells.Item(131, 11).Value = 4050  # K131: 3999 -> 4050
$ws.Cells.Item(131, 12).Value = 4494.4686  # L131: 4492.7628 -> 4494.4686
$ws.Cells.Item(131, 13).Value = 990  # M131: 1041 -> 990
$ws.Cells.Item(131, 14).Value = -14574.4686  # N131: -14572.7628 -> -14574.4686
$ws.Cells.Item(132, 8).Value = 1962.7142  # H132: 2092.25 -> 1962.7142
$ws.Cells.Item(132, 9).Value = 1962.7142  # I132: 2092.25 -> 1962.7142
$ws.Cells.Item(132, 11).Value = 17664.4278  # K132: 18830.25 -> 17664.4278
$ws.Cells.Item(132, 13).Value = -15134.4278  # M132: -16300.25 -> -15134.4278
$ws.Cells.Item(139, 8).Value = 17553.7  # H139: 16058.909 -> 17553.7
$ws.Cells.Item(139, 9).Value = 17553.7  # I139: 16058.909 -> 17553.7
$ws.Cells.Item(139, 11).Value = 52661.10000000001  # K139: 48176.727 -> 52661.10000000001
$ws.Cells.Item(139, 13).Value = -47521.10000000001  # M139: -43036.727 -> -47521.10000000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(36, 8).Value = 1924.3077  # H36: 2166.5833 -> 1924.3077
$ws.Cells.Item(36, 9).Value = 1668.4445  # I36: 1999.875 -> 1668.4445
$ws.Cells.Item(36, 11).Value = 1668.4445  # K36: 1999.875 -> 1668.4445
$ws.Cells.Item(36, 13).Value = -1183.4445  # M36: -1514.875 -> -1183.4445
$ws.Cells.Item(43, 8).Value = 4935.273  # H43: 6523.5 -> 4935.273
$ws.Cells.Item(43, 9).Value = 4935.273  # I43: 6523.5 -> 4935.273
$ws.Cells.Item(43, 11).Value = 4935.273  # K43: 6523.5 -> 4935.273
$ws.Cells.Item(43, 13).Value = -4784.273  # M43: -6372.5 -> -4784.273
$ws.Cells.Item(70, 8).Value = 17289.125  # H70: 15890.889 -> 17289.125
$ws.Cells.Item(70, 9).Value = 20390.273  # I70: 17976.154 -> 20390.273
$ws.Cells.Item(70, 10).Value = 10466.6  # J70: 10469.2 -> 10466.6
$ws.Cells.Item(70, 11).Value = 20390.273  # K70: 17976.154 -> 20390.273
$ws.Cells.Item(70, 12).Value = 10466.6  # L70: 10469.2 -> 10466.6
$ws.Cells.Item(70, 13).Value = -20120.273  # M70: -17706.154 -> -20120.273
$ws.Cells.Item(70, 14).Value = -11006.6  # N70: -11009.2 -> -11006.6
$ws.Cells.Item(73, 8).Value = 17289.125  # H73: 15890.889 -> 17289.125
$ws.Cells.Item(73, 9).Value = 20390.273  # I73: 17976.154 -> 20390.273
$ws.Cells.Item(73, 10).Value = 10466.6  # J73: 10469.2 -> 10466.6
$ws.Cells.Item(73, 11).Value = 20390.273  # K73: 17976.154 -> 20390.273
$ws.Cells.Item(73, 12).Value = 10466.6  # L73: 10469.2 -> 10466.6
$ws.Cells.Item(73, 13).Value = -19454.273  # M73: -17040.154 -> -19454.273
$ws.Cells.Item(73, 14).Value = -12338.6  # N73: -12341.2 -> -12338.6
$ws.Cells.Item(97, 8).Value = 3934.65  # H97: 3228.1924 -> 3934.65
$ws.Cells.Item(97, 9).Value = 2270.2307  # I97: 2000.7646 -> 2270.2307
$ws.Cells.Item(97, 10).Value = 7025.7144  # J97: 5546.6665 -> 7025.7144
$ws.Cells.Item(97, 11).Value = 2270.2307  # K97: 2000.7646 -> 2270.2307
$ws.Cells.Item(97, 12).Value = 7025.7144  # L97: 5546.6665 -> 7025.7144
$ws.Cells.Item(97, 13).Value = -1774.2307  # M97: -1504.7646 -> -1774.2307
$ws.Cells.Item(97, 14).Value = -8017.7144  # N97: -6538.6665 -> -8017.7144
$ws.Cells.Item(102, 8).Value = 3961  # H102: 4285.909 -> 3961
$ws.Cells.Item(102, 9).Value = 4124.476  # I102: 4517.8945 -> 4124.476
$ws.Cells.Item(102, 11).Value = 4124.476  # K102: 4517.8945 -> 4124.476
$ws.Cells.Item(102, 13).Value = -2502.476  # M102: -2895.8945 -> -2502.476
$ws.Cells.Item(107, 8).Value = 678.6667  # H107: 655.0454999999999 -> 678.6667
$ws.Cells.Item(107, 9).Value = 221.66667  # I107: 235.625 -> 221.66667
$ws.Cells.Item(107, 10).Value = 1021.4167  # J107: 894.7143 -> 1021.4167
$ws.Cells.Item(107, 11).Value = 221.66667  # K107: 235.625 -> 221.66667
$ws.Cells.Item(107, 12).Value = 1021.4167  # L107: 894.7143 -> 1021.4167
$ws.Cells.Item(107, 13).Value = 1698.33333  # M107: 1684.375 -> 1698.33333
$ws.Cells.Item(107, 14).Value = -4861.4167  # N107: -4734.7143 -> -4861.4167
$ws.Cells.Item(113, 8).Value = 3099.8235  # H113: 3231.125 -> 3099.8235
$ws.Cells.Item(113, 9).Value = 1683.0834  # I113: 1745.2727 -> 1683.0834
$ws.Cells.Item(113, 11).Value = 1683.0834  # K113: 1745.2727 -> 1683.0834
$ws.Cells.Item(113, 13).Value = 486.9166  # M113: 424.7273 -> 486.9166
$ws.Cells.Item(122, 8).Value = 2526.4  # H122: 2500.25 -> 2526.4
$ws.Cells.Item(122, 9).Value = 2526.4  # I122: 2571.7144 -> 2526.4
$ws.Cells.Item(122, 10).Value = 0  # J122: 2000 -> 0
$ws.Cells.Item(122, 11).Value = 7579.200000000001  # K122: 7715.1432 -> 7579.200000000001
$ws.Cells.Item(122, 12).Value = 0  # L122: 6000 -> 0
$ws.Cells.Item(122, 13).Value = -5129.200000000001  # M122: -5265.1432 -> -5129.200000000001
$ws.Cells.Item(122, 14).Value = $null  # N122: was -10900, now removed
$ws.Cells.Item(132, 8).Value = 14168.583  # H132: 13349.115 -> 14168.583
$ws.Cells.Item(132, 9).Value = 11096.762  # I132: 11098.238 -> 11096.762
$ws.Cells.Item(132, 10).Value = 35671.332  # J132: 22802.8 -> 35671.332
$ws.Cells.Item(132, 11).Value = 33290.286  # K132: 33294.714 -> 33290.286
$ws.Cells.Item(132, 12).Value = 107013.996  # L132: 68408.39999999999 -> 107013.996
$ws.Cells.Item(132, 13).Value = -30760.286  # M132: -30764.714 -> -30760.286
$ws.Cells.Item(132, 14).Value = -112073.996  # N132: -73468.39999999999 -> -112073.996
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 10488.5  # H7: 10289 -> 10488.5
$ws.Cells.Item(7, 9).Value = 12977  # I7: 12387.25 -> 12977
$ws.Cells.Item(7, 11).Value = 12977  # K7: 12387.25 -> 12977
$ws.Cells.Item(7, 13).Value = -12865  # M7: -12275.25 -> -12865
$ws.Cells.Item(16, 8).Value = 9998.5  # H16: 9999 -> 9998.5
$ws.Cells.Item(16, 10).Value = 9998  # J16: 0 -> 9998
$ws.Cells.Item(16, 12).Value = 9998  # L16: 0 -> 9998
$ws.Cells.Item(16, 14).Value = -10338  # N16: None -> -10338
$ws.Cells.Item(29, 8).Value = 9999  # H29: 0 -> 9999
$ws.Cells.Item(29, 9).Value = 9999  # I29: 0 -> 9999
$ws.Cells.Item(29, 11).Value = 9999  # K29: 0 -> 9999
$ws.Cells.Item(29, 13).Value = -9704  # M29: None -> -9704
$ws.Cells.Item(40, 8).Value = 6362.4116  # H40: 6556.9375 -> 6362.4116
$ws.Cells.Item(40, 10).Value = 7123.923  # J40: 7446.75 -> 7123.923
$ws.Cells.Item(40, 12).Value = 7123.923  # L40: 7446.75 -> 7123.923
$ws.Cells.Item(40, 14).Value = -7395.923  # N40: -7718.75 -> -7395.923
$ws.Cells.Item(46, 8).Value = 2409.125  # H46: 2474.3572 -> 2409.125
$ws.Cells.Item(46, 9).Value = 1283.5555  # I46: 1330.875 -> 1283.5555
$ws.Cells.Item(46, 10).Value = 3856.2856  # J46: 3999 -> 3856.2856
$ws.Cells.Item(46, 11).Value = 1283.5555  # K46: 1330.875 -> 1283.5555
$ws.Cells.Item(46, 12).Value = 3856.2856  # L46: 3999 -> 3856.2856
$ws.Cells.Item(46, 13).Value = -1095.5555  # M46: -1142.875 -> -1095.5555
$ws.Cells.Item(46, 14).Value = -4232.2856  # N46: -4375 -> -4232.2856
$ws.Cells.Item(55, 8).Value = 2170.7778  # H55: 2227.0857 -> 2170.7778
$ws.Cells.Item(55, 10).Value = 2948.65  # J55: 3093.3157 -> 2948.65
$ws.Cells.Item(55, 12).Value = 2948.65  # L55: 3093.3157 -> 2948.65
$ws.Cells.Item(55, 14).Value = -3294.65  # N55: -3439.3157 -> -3294.65
$ws.Cells.Item(93, 8).Value = 11182.25  # H93: 11692.947 -> 11182.25
$ws.Cells.Item(93, 10).Value = 7047.1  # J93: 7665.778 -> 7047.1
$ws.Cells.Item(93, 12).Value = 7047.1  # L93: 7665.778 -> 7047.1
$ws.Cells.Item(93, 14).Value = -9543.1  # N93: -10161.778 -> -9543.1
$ws.Cells.Item(122, 8).Value = 6275.75  # H122: 5381.0386 -> 6275.75
$ws.Cells.Item(122, 9).Value = 5847.3076  # I122: 4758.263 -> 5847.3076
$ws.Cells.Item(122, 11).Value = 17541.9228  # K122: 14274.789 -> 17541.9228
$ws.Cells.Item(122, 13).Value = -15091.9228  # M122: -11824.789 -> -15091.9228
$ws.Cells.Item(126, 8).Value = 10488.5  # H126: 10289 -> 10488.5
$ws.Cells.Item(126, 9).Value = 12977  # I126: 12387.25 -> 12977
$ws.Cells.Item(126, 11).Value = 38931  # K126: 37161.75 -> 38931
$ws.Cells.Item(126, 13).Value = -36461  # M126: -34691.75 -> -36461
$ws.Cells.Item(132, 8).Value = 857981.9399999999  # H132: 2119539.5 -> 857981.9399999999
$ws.Cells.Item(132, 9).Value = 2045.7307  # I132: 3397.25 -> 2045.7307
$ws.Cells.Item(132, 10).Value = 1917712.4  # J132: 2683844.2 -> 1917712.4
$ws.Cells.Item(132, 11).Value = 6137.1921  # K132: 10191.75 -> 6137.1921
$ws.Cells.Item(132, 12).Value = 5753137.199999999  # L132: 8051532.600000001 -> 5753137.199999999
$ws.Cells.Item(132, 13).Value = -3607.1921  # M132: -7661.75 -> -3607.1921
$ws.Cells.Item(132, 14).Value = -5758197.199999999  # N132: -8056592.600000001 -> -5758197.199999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 3195.9  # H2: 5266.5 -> 3195.9
$ws.Cells.Item(2, 9).Value = 3195.9  # I2: 5266.5 -> 3195.9
$ws.Cells.Item(2, 11).Value = 3195.9  # K2: 5266.5 -> 3195.9
$ws.Cells.Item(2, 13).Value = -3083.9  # M2: -5154.5 -> -3083.9
$ws.Cells.Item(110, 8).Value = 50000  # H110: 100644 -> 50000
$ws.Cells.Item(110, 10).Value = 50000  # J110: 100644 -> 50000
$ws.Cells.Item(110, 12).Value = 50000  # L110: 100644 -> 50000
$ws.Cells.Item(110, 14).Value = -58180  # N110: -108824 -> -58180
$ws.Cells.Item(113, 8).Value = 1850  # H113: 1804.4 -> 1850
$ws.Cells.Item(113, 9).Value = 2306.6155  # I113: 2230.5925 -> 2306.6155
$ws.Cells.Item(113, 11).Value = 6919.8465  # K113: 6691.7775 -> 6919.8465
$ws.Cells.Item(113, 13).Value = -4749.8465  # M113: -4521.7775 -> -4749.8465
$ws.Cells.Item(122, 8).Value = 2764.568  # H122: 2802.2327 -> 2764.568
$ws.Cells.Item(122, 9).Value = 1698.3871  # I122: 1711.2903 -> 1698.3871
$ws.Cells.Item(122, 10).Value = 5307  # J122: 5620.5 -> 5307
$ws.Cells.Item(122, 11).Value = 5095.1613  # K122: 5133.8709 -> 5095.1613
$ws.Cells.Item(122, 12).Value = 15921  # L122: 16861.5 -> 15921
$ws.Cells.Item(122, 13).Value = -2645.1613  # M122: -2683.8709 -> -2645.1613
$ws.Cells.Item(122, 14).Value = -20821  # N122: -21761.5 -> -20821
$ws.Cells.Item(126, 8).Value = 13381.6875  # H126: 9472.875 -> 13381.6875
$ws.Cells.Item(126, 9).Value = 10028.625  # I126: 5841.9375 -> 10028.625
$ws.Cells.Item(126, 11).Value = 30085.875  # K126: 17525.8125 -> 30085.875
$ws.Cells.Item(126, 13).Value = -27615.875  # M126: -15055.8125 -> -27615.875
$ws.Cells.Item(132, 8).Value = 5913.885  # H132: 6008.294 -> 5913.885
$ws.Cells.Item(132, 9).Value = 3054.9092  # I132: 3116.0312 -> 3054.9092
$ws.Cells.Item(132, 11).Value = 9164.7276  # K132: 9348.0936 -> 9164.7276
$ws.Cells.Item(132, 13).Value = -6634.7276  # M132: -6818.0936 -> -6634.7276
